$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two obsolete rows (neurobio-final, neuroelectronics) ---
$ws.Range("A5:D6").EntireRow.Delete()

# --- Row 2: nurb-discussion-1 / Neurobiology Discussion 1  ->  visual-coding / Basic Visual Coding ---
$ws.Range("A2").Value = "visual-coding"
$ws.Range("B2").Value = "Basic Visual Coding"

# --- Row 3: basic-visual-coding / Basic Visual Coding  ->  rsnn / Stochastic Computation... ---
$ws.Range("A3").Value = "rsnn"
$ws.Range("B3").Value = "Stochastic Computation in Recurrent Networks of Spiking Neurons"

# --- Row 4: nurb-discussion-2 / Neurobiology Discussion 2  ->  rsnn-slides / Stochastic Computation...(Slides) ---
$ws.Range("A4").Value = "rsnn-slides"
$ws.Range("B4").Value = "Stochastic Computation in Recurrent Networks of Spiking Neurons (Slides)"

# Rows 3 & 4 (cols A,B) lose their explicit "Arial 12" font override in the new file and
# fall back to the default/general style. Reproduce that by pasting the formatting of a
# never-formatted cell (outside the used range) onto them, which maps back onto the
# existing default style instead of fabricating a new one.
$ws.Range("E1").Copy()
$ws.Range("A3:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column B got wider to fit the new longer titles ---
$ws.Columns.Item(2).ColumnWidth = 67.58

# --- Update the (stale) _FilterDatabase defined name so it matches the shrunk data range ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "docs!_FilterDatabase") {
        $n.RefersTo = "=docs!`$D`$1:`$D`$2"
    }
}

# --- Selection moves to B2 ---
$null = $ws.Range("B2").Select()
